# Adds two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# header style used by the existing columns, and fills in data for rows 2-44.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style (bold, centered, bordered) used by the other header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data ---------------------------------------------------------------
$I0 = @(1,3,7,1,6,6,7,1,9,6,9,8,8,1,5,9,8,7,6,1,6,6,4,6,6,6,5,1,9,7,11,6,8,1,1,7,10,6,8,5,6,5,3)
$IF = @(2,5,7,1,7,7,7,3,10,7,9,9,9,2,6,9,9,8,7,2,8,7,6,8,7,8,6,1,9,8,11,7,8,2,2,7,10,8,8,5,7,6,3)

for ($r = 0; $r -lt $I0.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $I0[$r]
    $ws.Cells.Item($row, 10).Value = $IF[$r]
}

$ws.UsedRange | Out-Null
